$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.591.57'
$ws.Range('E2').Value = '  -0.23%  '
$ws.Range('D3').Value = '3.515.84'
$ws.Range('E3').Value = '  -1.98%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '''623.42'
$ws.Range('E5').Value = '  +4.33%  '
$ws.Range('D6').Value = '''172.89'
$ws.Range('E6').Value = '  -0.37%  '
$ws.Range('D7').Value = '''0.610'
$ws.Range('E7').Value = '  -1.34%  '
$ws.Range('D8').Value = '3.511.86'
$ws.Range('E8').Value = '  -1.91%  '
$ws.Range('E9').Value = '  -0.06%  '
$ws.Range('E10').Value = '  -0.80%  '
$ws.Range('D11').Value = '''7.16'
$ws.Range('E11').Value = '  -3.57%  '
$ws.Range('D12').Value = '''0.588'
$ws.Range('E12').Value = '  -0.43%  '
$ws.Range('D13').Value = '''46.46'
$ws.Range('E13').Value = '  -0.81%  '
$ws.Range('D14').Value = '''0.0000276'
$ws.Range('E14').Value = '  -0.91%  '
$ws.Range('D15').Value = '4.090.37'
$ws.Range('E15').Value = '  -1.92%  '
$ws.Range('E16').Value = '  -0.47%  '
$ws.Range('D17').Value = '''609.40'
$ws.Range('E17').Value = '  -0.69%  '
$ws.Range('D18').Value = '3.524.02'
$ws.Range('E18').Value = '  -1.52%  '
$ws.Range('D19').Value = '70.708.92'
$ws.Range('E19').Value = '  -0.20%  '
$ws.Range('D21').Value = '''17.74'
$ws.Range('E21').Value = '  +1.44%  '
$ws.Range('D22').Value = '''0.882'
$ws.Range('E22').Value = '  -0.56%  '
$ws.Range('D23').Value = '''9.11'
$ws.Range('E23').Value = '  -2.12%  '
$ws.Range('D24').Value = '''15.58'
$ws.Range('E24').Value = '  -2.39%  '
$ws.Range('D25').Value = '''97.33'
$ws.Range('E25').Value = '  +0.17%  '
$ws.Range('D26').Value = '''3.74'
$ws.Range('E26').Value = '  -1.13%  '
$ws.Range('E27').Value = '  +0.00%  '
$ws.Range('E28').Value = '  -2.97%  '
$ws.Range('D29').Value = '''33.55'
$ws.Range('E29').Value = '  -1.12%  '
$ws.Range('E30').Value = '  -1.39%  '
$ws.Range('D31').Value = '''3.02'
$ws.Range('E31').Value = '  -1.55%  '
$ws.Range('E32').Value = '  -4.27%  '
$ws.Range('D33').Value = '''1.29'
$ws.Range('E33').Value = '  -0.91%  '
$ws.Range('D34').Value = '''636.85'
$ws.Range('E34').Value = '  -0.89%  '
$ws.Range('D35').Value = '''6.83'
$ws.Range('E35').Value = '  -5.20%  '
$ws.Range('D36').Value = '''10.81'
$ws.Range('E36').Value = '  -0.37%  '
$ws.Range('D37').Value = '''0.0992'
$ws.Range('E37').Value = '  -2.14%  '
$ws.Range('D38').Value = '''0.0487'
$ws.Range('E38').Value = '  +1.09%  '
$ws.Range('D39').Value = '''3.42'
$ws.Range('E39').Value = '  -7.88%  '
$ws.Range('D40').Value = '''56.69'
$ws.Range('E40').Value = '  -0.98%  '
$ws.Range('E41').Value = '  +0.27%  '
$ws.Range('D42').Value = '''0.143'
$ws.Range('E42').Value = '  +0.36%  '
$ws.Range('D43').Value = '3.347.41'
$ws.Range('E43').Value = '  -1.61%  '
$ws.Range('D44').Value = '0.0₃0721'
$ws.Range('E44').Value = '  +0.38%  '
$ws.Range('E45').Value = '  +0.22%  '
$ws.Range('E46').Value = '  -3.81%  '
$ws.Range('D47').Value = '''31.98'
$ws.Range('E47').Value = '  -3.06%  '
$ws.Range('E48').Value = '  -5.16%  '
$ws.Range('D50').Value = '''132.99'
$ws.Range('E50').Value = '  +0.14%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '''0.156'
$ws.Range('E51').Value = '  +5.45%  '
